$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update row 2 (Target cluster becomes "ECs", and recompute derived metrics)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,4).Value  = "ECs"                  # D2 Target cluster
$ws.Cells.Item(2,5).Value  = 3                       # E2 Ligand-expressing cells
$ws.Cells.Item(2,6).Value  = 1                       # F2 Ligand detection rate
$ws.Cells.Item(2,7).Value  = 1.807599666666667       # G2 Ligand average expression value
$ws.Cells.Item(2,8).Value  = 5.422799                # H2 Ligand total expression value
$ws.Cells.Item(2,9).Value  = 1                       # I2 Ligand derived specificity (avg)
$ws.Cells.Item(2,10).Value = 1                       # J2 Ligand derived specificity (total)
$ws.Cells.Item(2,11).Value = 2                       # K2 Receptor-expressing cells
$ws.Cells.Item(2,12).Value = 0.6666666666666666      # L2 Receptor detection rate
$ws.Cells.Item(2,13).Value = 1.539665666666667       # M2 Receptor average expression value
$ws.Cells.Item(2,14).Value = 4.618997                # N2 Receptor total expression value
$ws.Cells.Item(2,15).Value = 0.3572088291809875      # O2 Receptor derived specificity (avg)
$ws.Cells.Item(2,16).Value = 0.3572088291809875      # P2 Receptor derived specificity (total)
$ws.Cells.Item(2,17).Value = 2.783099145844778        # Q2 Edge average expression weight
$ws.Cells.Item(2,18).Value = 25.047892312603          # R2 Edge total expression weight
$ws.Cells.Item(2,19).Value = 0.3572088291809875      # S2 Edge average expression derived specificity
$ws.Cells.Item(2,20).Value = 0.3572088291809875      # T2 Edge total expression derived specificity

# ---------------------------------------------------------------------------
# Update row 3 (Target cluster becomes "FAPs", and recompute derived metrics)
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,4).Value  = "FAPs"                  # D3 Target cluster
$ws.Cells.Item(3,5).Value  = 3                       # E3 Ligand-expressing cells
$ws.Cells.Item(3,6).Value  = 1                       # F3 Ligand detection rate
$ws.Cells.Item(3,7).Value  = 1.807599666666667       # G3 Ligand average expression value
$ws.Cells.Item(3,8).Value  = 5.422799                # H3 Ligand total expression value
$ws.Cells.Item(3,9).Value  = 1                       # I3 Ligand derived specificity (avg)
$ws.Cells.Item(3,10).Value = 1                       # J3 Ligand derived specificity (total)
$ws.Cells.Item(3,11).Value = 3                       # K3 Receptor-expressing cells
$ws.Cells.Item(3,12).Value = 1                       # L3 Receptor detection rate
$ws.Cells.Item(3,13).Value = 1.452872333333333       # M3 Receptor average expression value
$ws.Cells.Item(3,14).Value = 4.358617                # N3 Receptor total expression value
$ws.Cells.Item(3,15).Value = 0.3370724153789985      # O3 Receptor derived specificity (avg)
$ws.Cells.Item(3,16).Value = 0.3370724153789985      # P3 Receptor derived specificity (total)
$ws.Cells.Item(3,17).Value = 2.626211545442556        # Q3 Edge average expression weight
$ws.Cells.Item(3,18).Value = 23.635903908983          # R3 Edge total expression weight
$ws.Cells.Item(3,19).Value = 0.3370724153789985      # S3 Edge average expression derived specificity
$ws.Cells.Item(3,20).Value = 0.3370724153789985      # T3 Edge total expression derived specificity

# ---------------------------------------------------------------------------
# Add new row 4 (Target cluster "sCs")
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,1).Value  = "FAPs"                  # A4 Sending cluster
$ws.Cells.Item(4,2).Value  = "Gdf2"                  # B4 Ligand symbol
$ws.Cells.Item(4,3).Value  = "Acvr2b"                # C4 Receptor symbol
$ws.Cells.Item(4,4).Value  = "sCs"                   # D4 Target cluster
$ws.Cells.Item(4,5).Value  = 3                       # E4 Ligand-expressing cells
$ws.Cells.Item(4,6).Value  = 1                       # F4 Ligand detection rate
$ws.Cells.Item(4,7).Value  = 1.807599666666667       # G4 Ligand average expression value
$ws.Cells.Item(4,8).Value  = 5.422799                # H4 Ligand total expression value
$ws.Cells.Item(4,9).Value  = 1                       # I4 Ligand derived specificity (avg)
$ws.Cells.Item(4,10).Value = 1                       # J4 Ligand derived specificity (total)
$ws.Cells.Item(4,11).Value = 3                       # K4 Receptor-expressing cells
$ws.Cells.Item(4,12).Value = 1                       # L4 Receptor detection rate
$ws.Cells.Item(4,13).Value = 1.317729666666667       # M4 Receptor average expression value
$ws.Cells.Item(4,14).Value = 3.953189                # N4 Receptor total expression value
$ws.Cells.Item(4,15).Value = 0.3057187554400141      # O4 Receptor derived specificity (avg)
$ws.Cells.Item(4,16).Value = 0.3057187554400141      # P4 Receptor derived specificity (total)
$ws.Cells.Item(4,17).Value = 2.381927706223445        # Q4 Edge average expression weight
$ws.Cells.Item(4,18).Value = 21.437349356011          # R4 Edge total expression weight
$ws.Cells.Item(4,19).Value = 0.3057187554400141      # S4 Edge average expression derived specificity
$ws.Cells.Item(4,20).Value = 0.3057187554400141      # T4 Edge total expression derived specificity
